$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "'-335"
$ws.Range("E3").Value = "'+260"

$ws.Range("D4").Value = "'+190"
$ws.Range("E4").Value = "'-235"

$ws.Range("D5").Value = "'-210"
$ws.Range("E5").Value = "'+175"

$ws.Range("D6").Value = "'+135"
$ws.Range("E6").Value = "'-155"

$ws.Range("D8").Value = "'-265"
$ws.Range("E8").Value = "'+215"

$ws.Range("D9").Value = "'+155"
$ws.Range("E9").Value = "'-180"

$ws.Range("D11").Value = "'-145"
$ws.Range("E11").Value = "'+125"

$ws.Range("D12").Value = "'+125"
$ws.Range("E12").Value = "'-145"

$ws.Range("D13").Value = "'-105"
$ws.Range("E13").Value = "'-115"

$ws.Range("D15").Value = "'-190"
$ws.Range("E15").Value = "'+160"
